$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QA")

# --- Update the test-case identifiers / tester name -----------------------
# Row 3: the tester's first name changes from "Axel" to "Raúl"
$ws.Range("B3").Value = "Raúl"

# Row 2: DNI-like code changes (PRU-888-011 -> AUT_JF_QA_001)
$ws.Range("A2").Value = "AUT_JF_QA_001"

# Row 3: DNI-like code changes (PRU-888-012 -> AUT_TF_QA_001)
$ws.Range("A3").Value = "AUT_TF_QA_001"

# Row 4: a new code is added (the row previously had no DNI value)
$ws.Range("A4").Value = "AUT_TF_QA_002"

# The now orphaned, style-only cell A6 is removed entirely
$ws.Range("A6").Clear()

# --- Resize column A (the user manually widened it, dropping autofit) -----
$ws.Columns.Item(1).ColumnWidth = 17.5
